# Beginning of Second Question of Assignment
# Append the new "Socks In The Dark" question block after the end of the
# first question, right before the trailing bookmark ("_GoBack") paragraph
# that closes the document body.  The bookmark paragraph itself becomes the
# paragraph that holds the new "It is dark and I can only pick..." text, so
# the bookmark stays exactly where it was (at the very end of the body).

$d = $word.ActiveDocument

# Locate the document's final paragraph -- it is the (originally empty)
# paragraph that carries the _GoBack bookmark.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

# A collapsed range positioned right at the start of that paragraph: any
# OOXML we insert there lands *before* the bookmark, and a trailing partial
# paragraph in the inserted fragment merges into this same paragraph,
# keeping the bookmark as the last thing in the body.
$insertionPoint = $d.Range($last.Range.Start, $last.Range.Start)

$newParagraphsXml = @'
<w:p/><w:p/><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:i/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Socks In The Dark:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:i/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Problem: </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>It is dark and I can only pick one pair of socks out of 20 socks that are of 3 different colors. The only way I can veify</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newParagraphsXml)
